$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1147919382644886
$ws.Range("D2").Value = 0.04883252629998935
$ws.Range("E2").Value = 0.07691568776725166
$ws.Range("F2").Value = 1.68497475861966
$ws.Range("G2").Value = 1.651992108811015
$ws.Range("H2").Value = 1.329248330178615
$ws.Range("I2").Value = 1.501420618137686
$ws.Range("L2").Value = 0.09724694543284329
$ws.Range("M2").Value = 1.598258328719723
$ws.Range("N2").Value = 1.647826646260853
$ws.Range("C3").Value = 0.1154208897993367
$ws.Range("D3").Value = 0.04901561773452556
$ws.Range("E3").Value = 0.0777627994897987
$ws.Range("F3").Value = 1.62426841647131
$ws.Range("G3").Value = 1.58111630897281
$ws.Range("H3").Value = 1.304050679731603
$ws.Range("I3").Value = 1.461805080788452
$ws.Range("L3").Value = 0.09860868902240227
$ws.Range("M3").Value = 1.449215974917067
$ws.Range("N3").Value = 1.527592373724246
$ws.Range("C4").Value = 0.1158309581009576
$ws.Range("D4").Value = 0.04916958880606614
$ws.Range("E4").Value = 0.07831448098534555
$ws.Range("F4").Value = 1.588355854453297
$ws.Range("G4").Value = 1.538998674695506
$ws.Range("H4").Value = 1.289562749411658
$ws.Range("I4").Value = 1.438700112296686
$ws.Range("L4").Value = 0.09949819799485482
$ws.Range("M4").Value = 1.357760496078129
$ws.Range("N4").Value = 1.454045031664094
$ws.Range("C5").Value = 0.1160040861500349
$ws.Range("D5").Value = 0.04924269305971407
$ws.Range("E5").Value = 0.07854724125184553
$ws.Range("F5").Value = 1.574059940776536
$ws.Range("G5").Value = 1.522183227385028
$ws.Range("H5").Value = 1.283904276851018
$ws.Range("I5").Value = 1.42958834623802
$ws.Range("L5").Value = 0.09987411229118415
$ws.Range("M5").Value = 1.320507426216395
$ws.Range("N5").Value = 1.424146867209203
$ws.Range("C6").Value = 0.1160331980671785
$ws.Range("D6").Value = 0.04925545509309615
$ws.Range("E6").Value = 0.07858637115352174
$ws.Range("F6").Value = 1.571706470798588
$ws.Range("G6").Value = 1.519411918256054
$ws.Range("H6").Value = 1.282979464675151
$ws.Range("I6").Value = 1.428093603855018
$ws.Range("L6").Value = 0.09993734402684495
$ws.Range("M6").Value = 1.314322578023763
$ws.Range("N6").Value = 1.419186806126646
$ws.Range("C7").Value = 0.1158332685650532
$ws.Range("D7").Value = 0.04917053288589202
$ws.Range("E7").Value = 0.07831758788452792
$ws.Range("F7").Value = 1.588161687883357
$ws.Range("G7").Value = 1.538770492991603
$ws.Range("H7").Value = 1.289485445752206
$ws.Range("I7").Value = 1.438576001931267
$ws.Range("L7").Value = 0.09950321331122858
$ws.Range("M7").Value = 1.3572580223492
$ws.Range("N7").Value = 1.453641514640168
$ws.Range("C8").Value = 0.1150038513541904
$ws.Range("D8").Value = 0.04888697817437304
$ws.Range("E8").Value = 0.07720123047155347
$ws.Range("F8").Value = 1.6637584667087
$ws.Range("G8").Value = 1.627260848247602
$ws.Range("H8").Value = 1.320354923021426
$ws.Range("I8").Value = 1.487506400153521
$ws.Range("L8").Value = 0.09770539660343047
$ws.Range("M8").Value = 1.546857457174696
$ws.Range("N8").Value = 1.606314316817191
$ws.Range("C9").Value = 0.1135662494268352
$ws.Range("D9").Value = 0.04866471923077853
$ws.Range("E9").Value = 0.07526186120419265
$ws.Range("F9").Value = 1.822978065962076
$ws.Range("G9").Value = 1.81210829255761
$ws.Range("H9").Value = 1.388782503755863
$ws.Range("I9").Value = 1.593269247417794
$ws.Range("L9").Value = 0.09460324308338031
$ws.Range("M9").Value = 1.919073313275106
$ws.Range("N9").Value = 1.907774923655353
$ws.Range("C10").Value = 0.1126242532129318
$ws.Range("D10").Value = 0.0487107069137025
$ws.Range("E10").Value = 0.07398852252564936
$ws.Range("F10").Value = 1.946902114985164
$ws.Range("G10").Value = 1.955119412023066
$ws.Range("H10").Value = 1.44399912892635
$ws.Range("I10").Value = 1.677155596836329
$ws.Range("L10").Value = 0.09258176220916425
$ws.Range("M10").Value = 2.192766726939936
$ws.Range("N10").Value = 2.130374038076184
$ws.Range("C11").Value = 0.1122203132440696
$ws.Range("D11").Value = 0.04877835264119312
$ws.Range("E11").Value = 0.07344198368700194
$ws.Range("F11").Value = 2.004843344694848
$ws.Range("G11").Value = 2.021810678652628
$ws.Range("H11").Value = 1.470221249066128
$ws.Range("I11").Value = 1.716705012058
$ws.Range("L11").Value = 0.09171800614218384
$ws.Range("M11").Value = 2.317323614499514
$ws.Range("N11").Value = 2.231854144824979
$ws.Range("C12").Value = 0.1120708706129356
$ws.Range("D12").Value = 0.04881079352490758
$ws.Range("E12").Value = 0.07323971528241202
$ws.Range("F12").Value = 2.027014023173791
$ws.Range("G12").Value = 2.047305290251529
$ws.Range("H12").Value = 1.480311908192675
$ws.Range("I12").Value = 1.731884621960319
$ws.Range("L12").Value = 0.09139894781915459
$ws.Range("M12").Value = 2.364496824869036
$ws.Range("N12").Value = 2.270310682055594
$ws.Range("C13").Value = 0.1121028993935269
$ws.Range("D13").Value = 0.0488035015550139
$ws.Range("E13").Value = 0.07328306878711288
$ws.Range("F13").Value = 2.022228892367878
$ws.Range("G13").Value = 2.041803811591365
$ws.Range("H13").Value = 1.47813151236835
$ws.Range("I13").Value = 1.728606334191184
$ws.Range("L13").Value = 0.09146730571890416
$ws.Range("M13").Value = 2.354336973621344
$ws.Range("N13").Value = 2.262027174025036
$ws.Range("C14").Value = 0.1122079480283737
$ws.Range("D14").Value = 0.04878088417741111
$ws.Range("E14").Value = 0.07342524892373703
$ws.Range("F14").Value = 2.006662714188423
$ws.Range("G14").Value = 2.023903295294986
$ws.Range("H14").Value = 1.471048176096076
$ws.Range("I14").Value = 1.717949757759328
$ws.Range("L14").Value = 0.09169159614224043
$ws.Range("M14").Value = 2.321204465559333
$ws.Range("N14").Value = 2.235017439863952
$ws.Range("C15").Value = 0.112272751437434
$ws.Range("D15").Value = 0.04876792236027683
$ws.Range("E15").Value = 0.07351294937612884
$ws.Range("F15").Value = 1.997158005556258
$ws.Range("G15").Value = 2.012970126900484
$ws.Range("H15").Value = 1.466730449086242
$ws.Range("I15").Value = 1.711448847766306
$ws.Range("L15").Value = 0.09183002589763589
$ws.Range("M15").Value = 2.300910619772424
$ws.Range("N15").Value = 2.218476778851027
$ws.Range("C16").Value = 0.1126511450543468
$ws.Range("D16").Value = 0.04870723534868659
$ws.Range("E16").Value = 0.07402489749944774
$ws.Range("F16").Value = 1.943147418249026
$ws.Range("G16").Value = 1.950794302208067
$ws.Range("H16").Value = 1.442307848937133
$ws.Range("I16").Value = 1.674599187116939
$ws.Range("L16").Value = 0.09263933402030133
$ws.Range("M16").Value = 2.184627582015537
$ws.Range("N16").Value = 2.123746197686387
$ws.Range("C17").Value = 0.1128895627356705
$ws.Range("D17").Value = 0.04868204810234289
$ws.Range("E17").Value = 0.07434733189870979
$ws.Range("F17").Value = 1.910418112484905
$ws.Range("G17").Value = 1.913073685006452
$ws.Range("H17").Value = 1.427609627131375
$ws.Range("I17").Value = 1.652351262939661
$ws.Range("L17").Value = 0.09315011629501768
$ws.Range("M17").Value = 2.11330416961664
$ws.Range("N17").Value = 2.06568574622986
$ws.Range("C18").Value = 0.1130290085779144
$ws.Range("D18").Value = 0.04867195266700719
$ws.Range("E18").Value = 0.07453586696080272
$ws.Range("F18").Value = 1.891740373224707
$ws.Range("G18").Value = 1.891531354255164
$ws.Range("H18").Value = 1.419259360446802
$ws.Range("I18").Value = 1.639685389457213
$ws.Range("L18").Value = 0.09344915932705078
$ws.Range("M18").Value = 2.07228583452661
$ws.Range("N18").Value = 2.032311679205634
$ws.Range("C19").Value = 0.1130766204449394
$ws.Range("D19").Value = 0.04866928559422945
$ws.Range("E19").Value = 0.07460023086224599
$ws.Range("F19").Value = 1.88544158295656
$ws.Range("G19").Value = 1.8842637162399
$ws.Range("H19").Value = 1.416449858837098
$ws.Range("I19").Value = 1.635419262099887
$ws.Range("L19").Value = 0.09355131265295924
$ws.Range("M19").Value = 2.058398631615177
$ws.Range("N19").Value = 2.021015459778539
$ws.Range("C20").Value = 0.1128639433442657
$ws.Range("D20").Value = 0.04868427405266118
$ws.Range("E20").Value = 0.07431268957784276
$ws.Range("F20").Value = 1.91388692877905
$ws.Range("G20").Value = 1.917073179210746
$ws.Range("H20").Value = 1.429163523285894
$ws.Range("I20").Value = 1.654706057519689
$ws.Range("L20").Value = 0.09309519882443951
$ws.Range("M20").Value = 2.120896160276033
$ws.Range("N20").Value = 2.07186425440824
$ws.Range("C21").Value = 0.1121769972573503
$ws.Range("D21").Value = 0.04878734137103891
$ws.Range("E21").Value = 0.07338335985342148
$ws.Range("F21").Value = 2.011228613310976
$ws.Range("G21").Value = 2.029154555017953
$ws.Range("H21").Value = 1.473124338945297
$ws.Range("I21").Value = 1.721074313349362
$ws.Range("L21").Value = 0.09162549874846349
$ws.Range("M21").Value = 2.330936127682065
$ws.Range("N21").Value = 2.242950112879271
$ws.Range("C22").Value = 0.1117485530589661
$ws.Range("D22").Value = 0.04889453311061942
$ws.Range("E22").Value = 0.07280334494511642
$ws.Range("F22").Value = 2.076187262530794
$ws.Range("G22").Value = 2.103807940021284
$ws.Range("H22").Value = 1.502794231680468
$ws.Range("I22").Value = 1.765635253134747
$ws.Range("L22").Value = 0.09071175447650859
$ws.Range("M22").Value = 2.46824507622361
$ws.Range("N22").Value = 2.35492766920936
$ws.Range("C23").Value = 0.11197534922443
$ws.Range("D23").Value = 0.04883364283845992
$ws.Range("E23").Value = 0.07311040990202677
$ws.Range("F23").Value = 2.041393577744088
$ws.Range("G23").Value = 2.06383409301867
$ws.Range("H23").Value = 1.486872203662472
$ws.Range("I23").Value = 1.741742648787437
$ws.Range("L23").Value = 0.09119515573366499
$ws.Range("M23").Value = 2.394957859694216
$ws.Range("N23").Value = 2.295149267256704
$ws.Range("C24").Value = 0.1128755184691101
$ws.Range("D24").Value = 0.04868325404933671
$ws.Range("E24").Value = 0.07432834151753021
$ws.Range("F24").Value = 1.912318244995731
$ws.Range("G24").Value = 1.915264560494393
$ws.Range("H24").Value = 1.428460695718684
$ws.Range("I24").Value = 1.653641066560681
$ws.Range("L24").Value = 0.0931200102323082
$ws.Range("M24").Value = 2.117463863419545
$ws.Range("N24").Value = 2.069070932986278
$ws.Range("C25").Value = 0.1139350345970129
$ws.Range("D25").Value = 0.04868851898792315
$ws.Range("E25").Value = 0.07575984874644881
$ws.Range("F25").Value = 1.778703923253801
$ws.Range("G25").Value = 1.760859342914756
$ws.Range("H25").Value = 1.369412631883733
$ws.Range("I25").Value = 1.563586662540928
$ws.Range("L25").Value = 0.0953971791133128
$ws.Range("M25").Value = 1.818337707854283
$ws.Range("N25").Value = 1.826015578159428
